# Generate Report for Handoff
# Update status text + timestamps, and shrink the "Latest Handoff/Handback"
# related date columns that previously accommodated the longer status text
# (was ~29.98 characters wide, now ~17.22 characters wide -> ColumnWidth 16.3
# is the input that lands on the nearest representable column width).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-20 09:29:19"

$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-20 09:29:08"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-10-20 09:29:19"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
